$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "[J. Paulo S.-T. NãoMetalicos-1NA, J. Paulo S.-T. NãoMetalicos-1NA, J. Paulo S.-T. NãoMetalicos-1NA, J. Paulo S.-T. NãoMetalicos-1NA]"
$ws.Range("C18").Value = "Tiago P.-M.T.F.-"
$ws.Range("D18").Value = "[Suzanny-Metalografia-1NA, Suzanny-Metalografia-1NA, Suzanny-Metalografia-1NA, Suzanny-Metalografia-1NA]"
$ws.Range("E18").Value = "[Anderson-Metrologia 1-1NA, Anderson-Metrologia 1-1NA, Anderson-Metrologia 1-1NA, Anderson-Metrologia 1-1NA]"
$ws.Range("F18").Value = "[Emerson-Comandos Eletricos-1NA, Emerson-Comandos Eletricos-1NA, Emerson-Comandos Eletricos-1NA, Emerson-Comandos Eletricos-1NA]"

# Row 19
$ws.Range("B19").Value = "[Anderson-Tornearia-1NA, Anderson-Tornearia-1NA, Anderson-Tornearia-1NA, Anderson-Tornearia-1NA]"
$ws.Range("C19").Value = "[Joel L.-T. M. Metalicos-1NA, Joel L.-T. M. Metalicos-1NA, Joel L.-T. M. Metalicos-1NA, Joel L.-T. M. Metalicos-1NA]"
$ws.Range("E19").Value = "[Humberto-Desenho tecnico mecanico – T1-1NA, Humberto-Desenho tecnico mecanico – T1-1NA, Humberto-Desenho tecnico mecanico – T1-1NA, Humberto-Desenho tecnico mecanico – T1-1NA]"

# Row 20
$ws.Range("B20").Value = "[Aselmo-Manut. Mot. End.-1NA, Aselmo-Manut. Mot. End.-1NA, Aselmo-Manut. Mot. End.-1NA, Aselmo-Manut. Mot. End.-1NA]"
$ws.Range("C20").Value = "[Rachel-Trat. Termicos-1NA, Rachel-Trat. Termicos-1NA, Rachel-Trat. Termicos-1NA, Rachel-Trat. Termicos-1NA]"
$ws.Range("D20").Value = "[Elcio D.-Desenho tecnico mecanico – T2-1NA, Elcio D.-Desenho tecnico mecanico – T2-1NA, Elcio D.-Desenho tecnico mecanico – T2-1NA, Elcio D.-Desenho tecnico mecanico – T2-1NA]"
$ws.Range("E20").Value = "[Gisele-E. D. N. D.-1NA, Gisele-E. D. N. D.-1NA, Gisele-E. D. N. D.-1NA, Gisele-E. D. N. D.-1NA]"
$ws.Range("F20").Value = "Gilberto-M.T.R.M.-"

# Row 21
$ws.Range("B21").Value = "[Victor S.-Ajustagem-1NA, Victor S.-Ajustagem-1NA, Victor S.-Ajustagem-1NA, Victor S.-Ajustagem-1NA]"
$ws.Range("C21").Value = "[Valmir-Caldeiraria-1NA, Valmir-Caldeiraria-1NA, Valmir-Caldeiraria-1NA, Valmir-Caldeiraria-1NA]"
$ws.Range("D21").Value = "Tiago P.-M.T.F.-"
$ws.Range("F21").Value = "[Clesidson-Elet. Digi. Básica-1NA, Clesidson-Elet. Digi. Básica-1NA, Clesidson-Elet. Digi. Básica-1NA, Clesidson-Elet. Digi. Básica-1NA]"
